$d = $word.ActiveDocument

# 1) Collapse the multi-run, proofErr-laden "Mein Leben..." sentence back into
#    a single plain run. Word's Find/Replace naturally merges the runs it
#    rewrites and drops the now-irrelevant grammar/spelling proofing marks.
$oldSentence = "Mein Leben dreht sich um Lernen, Arbeiten(man muss Geld verdienen) und Schlafen. Eine franzözische Redewendung dafür ist:Metro-Boulot-Dodo."
$d.Content.Find.Execute($oldSentence, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $oldSentence, 2) | Out-Null

# 2) Append the closing: a blank paragraph followed by the sign-off line,
#    both using the same "NurText" / bold-italic Times New Roman run
#    formatting as the rest of the letter.
$end = $d.Content.End
$insertionRange = $d.Range($end, $end)

$closingXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="NurText"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:bCs/>
                <w:i/>
                <w:iCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="NurText"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:bCs/>
                <w:i/>
                <w:iCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:bCs/>
                <w:i/>
                <w:iCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>Mit freundlichen Gr&#252;&#223;en</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionRange.InsertXML($closingXml) | Out-Null
